# Update "想去人数" (F column) counts across sheets to match the latest scrape.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1550
$ws1.Range("F5").Value  = 8443
$ws1.Range("F10").Value = 228
$ws1.Range("F14").Value = 266
$ws1.Range("F17").Value = 1367
$ws1.Range("F18").Value = 1296
$ws1.Range("F23").Value = 194
$ws1.Range("F26").Value = 52
$ws1.Range("F27").Value = 257
$ws1.Range("F31").Value = 190
$ws1.Range("F37").Value = 108
$ws1.Range("F41").Value = 1215

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 120
$ws2.Range("F35").Value = 153

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F9").Value  = 1902
$ws3.Range("F10").Value = 2879

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1550
$ws4.Range("F7").Value  = 8443
$ws4.Range("F11").Value = 1902
$ws4.Range("F12").Value = 2879
$ws4.Range("F17").Value = 228
$ws4.Range("F20").Value = 266
$ws4.Range("F21").Value = 1367
$ws4.Range("F22").Value = 1296
$ws4.Range("F26").Value = 194
$ws4.Range("F33").Value = 190
